# Append a new row (row 5) to Sheet1, mirroring the existing rows 2-4,
# with the timestamp "٠٥‏/٠٥‏/٢٠٢٥ ٠٢:١٠:٢٩ م" (05/05/2025 02:10:29 PM).
#
# Columns A-H: ملاحظات, المرافق, الكمية, المخيم, نوع المسافة, المركبة, المؤسسة, الوقت
# Row 5 mirrors rows 2-4 except for the notes (A, blank) and the time (H).
#
# A5 and C5 are written with a leading apostrophe so the engine stores them
# as literal text (A5 as an explicit empty string, C5 as the text "233",
# matching how the source file already stores A2:A4/C2:C4 as text) instead
# of coercing "233" to a number or dropping the empty value entirely. The
# Style reset afterwards clears the transient "quote prefix" formatting so
# the cells keep the workbook's default (unstyled) appearance, just like
# every other cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "'"
$ws.Range("B5").Value = "احمد"
$ws.Range("C5").Value = "'233"
$ws.Range("D5").Value = "الصمود"
$ws.Range("E5").Value = "الرحلة 2"
$ws.Range("F5").Value = "C2"
$ws.Range("G5").Value = "IDRF"
$ws.Range("H5").Value = "٠٥‏/٠٥‏/٢٠٢٥ ٠٢:١٠:٢٩ م"

$ws.Range("A5").Style = "Normal"
$ws.Range("C5").Style = "Normal"
